# Round 2 dense-answer relevance scoring:
# - add an "Answer_relevance" 0/1 score in column D for every data row
#   (rows 2-71 already/newly), rows whose B value is "tell me about how
#   RAG works." (every 14th row starting at 15) score 0, all others score 1
# - column B gets re-styled: narrower + word-wrap (matches the wrap style
#   already used by column C, style index 3 / header style index 2)
# - selection / zoom / scroll position updated to match where the author
#   was last working (bottom of the sheet, zoomed in)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B: narrower width + word wrap (reuses existing wrap styles) ---
$ws.Range("B1").WrapText = $true
$ws.Columns("B").ColumnWidth = 20.5

for ($r = 2; $r -le 71; $r++) {
    $ws.Cells.Item($r, 2).WrapText = $true
}

# --- Column D: Answer_relevance scores for rows 7-71 (2-6 already set) ---
for ($r = 7; $r -le 71; $r++) {
    $cycleRow = ($r - 1) % 14
    if ($cycleRow -eq 0) {
        $ws.Cells.Item($r, 4).Value = 0
    } else {
        $ws.Cells.Item($r, 4).Value = 1
    }
}

# --- View state: scroll near the bottom of the sheet, zoomed to 120% ---
$ws.Activate()
try { $excel.ActiveWindow.ScrollRow = 69 } catch {}
try { $excel.ActiveWindow.ScrollColumn = 2 } catch {}
try { $excel.ActiveWindow.Zoom = 120 } catch {}

$ws.Range("D15").Select()
